# [Fonds de solidarite] Add 2022-06-09 data
# Updates "nombre_aides" (col C) and "montant_total" (col E) for the rows
# whose figures changed with the new data extraction.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 7;   C = 7014;   E = 290682013 },
    @{ Row = 37;  C = 23050;  E = 130197201 },
    @{ Row = 51;  C = 6357;   E = 12201689 },
    @{ Row = 92;  C = 409282; E = 1597152985 },
    @{ Row = 93;  C = 209656; E = 1310025277 },
    @{ Row = 94;  C = 94231;  E = 919017349 },
    @{ Row = 95;  C = 50802;  E = 934242292 },
    @{ Row = 116; C = 4566;   E = 20667158 },
    @{ Row = 121; C = 14;     E = 1153896 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
